$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Question A block (GCI / Richardson extrapolation convergence check).
# Text labels are poked in the same (slightly non-linear) order the author's
# shared-string table was built, so new-string allocation order matches; the
# numeric/formula siblings are filled in right alongside each label since
# their write order has no bearing on the shared-string table.
# ---------------------------------------------------------------------------
$ws.Range("A29").Value = "p obs"
$ws.Range("B29").Formula = "=(LN((E25^2-1)*(C23-C24)/(C24-C25)+E25^2))/LN(E25*E24)"

$ws.Range("A28").Value = "p formel"
$ws.Range("B28").Value = 2

$ws.Range("C30").Value = "<10%"

$ws.Range("A31").Value = "GCI"
$ws.Range("B31").Formula = "=1.25/(E25^B28-1)*ABS(C24-C25)"

$ws.Range("A30").Value = "critère GCI"
$ws.Range("B30").Formula = "=ABS((B29-B28)/B28)"
$ws.Range("B30").Style = "Percent"
$ws.Range("B30").NumberFormat = "0.0000%"
$ws.Range("B30").Font.Name = "Aptos Narrow"
$ws.Range("B30").Font.Size = 11
$ws.Range("B30").Font.ThemeColor = 1

$ws.Range("C31").Value = "microns"

$ws.Range("A32").Value = "fh"
$ws.Range("B32").Formula = "=C25"
$ws.Range("C32").Value = "microns"
$ws.Range("D32").NumberFormat = "0.00E+00"

$ws.Range("A33").Value = "u_num"
$ws.Range("B33").Formula = "=B31/2"
$ws.Range("C33").Value = "microns"

$ws.Range("A27").Value = "QUESTION A"

# ---------------------------------------------------------------------------
# Question D block, part 1 (combined numerical uncertainty)
# ---------------------------------------------------------------------------
$ws.Range("A36").Value = "u_input"
$ws.Range("B36").Value = 2.685

$ws.Range("A35").Value = "QUESTION D"

$ws.Range("A37").Value = "u_D"
$ws.Range("B37").Value = 17.7789

$ws.Range("A38").Value = "u_val"
$ws.Range("B38").Formula = "=SQRT(B33^2+B36^2+B37^2)"

# ---------------------------------------------------------------------------
# Question D block, part 2 (validation interval E +/- k*u_val)
# ---------------------------------------------------------------------------
$ws.Range("A40").Value = "E"
$ws.Range("B40").Value = -58.605

$ws.Range("A41").Value = "k"
$ws.Range("B41").Value = 2

$ws.Range("A42").Value = "E-ku_val"
$ws.Range("B42").Formula = "=B40-B41*B38"
$ws.Range("C42").Formula = "=B41*B38"

$ws.Range("A43").Value = "E+ku_val"
$ws.Range("B43").Formula = "=B40+B41*B38"

# ---------------------------------------------------------------------------
# View bookkeeping to mirror the final selection / scroll position
# ---------------------------------------------------------------------------
$ws.Range("C43").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 21
